# Update the "Spots" template dates so the whole example campaign falls
# inside a single month (per commit: "aktualizace vzoru, aby datumy byly
# v jednom mesici" -> update the template so the dates are in one month).
#
# The dates are stored as plain date-serial numbers in column F (Start
# date) and column G (End date), formatted with a date number format
# (style s="3"), so we only need to update the underlying numeric value -
# Excel keeps the existing date formatting on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spots")

$ws.Range("F3").Value  = 45127
$ws.Range("F4").Value  = 45128
$ws.Range("F5").Value  = 45129
$ws.Range("F6").Value  = 45127
$ws.Range("F7").Value  = 45128
$ws.Range("F8").Value  = 45130
$ws.Range("F9").Value  = 45126
$ws.Range("F10").Value = 45127
$ws.Range("F11").Value = 45129
$ws.Range("F12").Value = 45119
$ws.Range("G12").Value = 45127
$ws.Range("F13").Value = 45110
$ws.Range("G13").Value = 45110
$ws.Range("F14").Value = 45127
$ws.Range("G14").Value = 45129
